$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (some look numeric, e.g. "0.7405");
# force text format so Excel does not auto-convert them to numbers,
# matching the original inlineStr cell type, then restore the default style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.952.62'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '1.876.55'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '0.7405'
$ws.Range('E5').Value = '  -3.79%  '
$ws.Range('D6').Value = '242.69'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '0.3148'
$ws.Range('E8').Value = '  +0.97%  '
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('D10').Value = '24.64'
$ws.Range('E10').Value = '  -3.71%  '
$ws.Range('D11').Value = '0.08330'
$ws.Range('E11').Value = '  -3.24%  '
$ws.Range('D12').Value = '0.7524'
$ws.Range('E12').Value = '  -1.35%  '
$ws.Range('D13').Value = '1.893.01'
$ws.Range('E13').Value = '  -1.55%  '
$ws.Range('D14').Value = '5.421'
$ws.Range('E14').Value = '  +0.99%  '
$ws.Range('D15').Value = '92.65'
$ws.Range('E15').Value = '  -0.98%  '
$ws.Range('D16').Value = '29.980.34'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('D17').Value = '6.121'
$ws.Range('E17').Value = '  -0.64%  '
$ws.Range('D18').Value = '250.01'
$ws.Range('E18').Value = '  +2.28%  '
$ws.Range('D19').Value = '13.58'
$ws.Range('E19').Value = '  -1.30%  '
$ws.Range('D20').Value = '0.000007859'
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').Value = '2.142.46'
$ws.Range('E22').Value = '  -1.91%  '
$ws.Range('D23').Value = '8.041'
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').Value = '0.1552'
$ws.Range('E25').Value = '  -5.77%  '
$ws.Range('D26').Value = '9.266'
$ws.Range('E26').Value = '  -1.04%  '
$ws.Range('D27').Value = '165.09'
$ws.Range('E27').Value = '  +1.73%  '
$ws.Range('D28').Value = '18.70'
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('D29').Value = '2.036'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').Value = '1.516'
$ws.Range('E30').Value = '  +4.30%  '
$ws.Range('D31').Value = '4.612'
$ws.Range('E31').Value = '  +2.49%  '
$ws.Range('D32').Value = '1.531'
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('D33').Value = '4.296'
$ws.Range('E33').Value = '  +4.79%  '
$ws.Range('D34').Value = '0.05325'
$ws.Range('E34').Value = '  -2.25%  '
$ws.Range('D35').Value = '1.236'
$ws.Range('E35').Value = '  -0.26%  '
$ws.Range('D36').Value = '0.7483'
$ws.Range('E36').Value = '  +0.51%  '
$ws.Range('D37').Value = '1.002'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').Value = '2.700'
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').Value = '0.01970'
$ws.Range('E39').Value = '  +0.67%  '
$ws.Range('D40').Value = '2.758'
$ws.Range('E40').Value = '  -0.81%  '
$ws.Range('D41').Value = '0.4562'
$ws.Range('E41').Value = '  +2.18%  '
$ws.Range('D42').Value = '1.112.92'
$ws.Range('E42').Value = '  +0.41%  '
$ws.Range('D43').Value = '6.150'
$ws.Range('E43').Value = '  +1.33%  '
$ws.Range('D44').Value = '72.36'
$ws.Range('E44').Value = '  -0.90%  '
$ws.Range('D45').Value = '0.8565'
$ws.Range('E45').Value = '  +0.78%  '
$ws.Range('D46').Value = '1.002'
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range('D47').Value = '104.08'
$ws.Range('E47').Value = '  +1.74%  '
$ws.Range('D48').Value = '1.857'
$ws.Range('E48').Value = '  -0.30%  '
$ws.Range('D49').Value = '7.613'
$ws.Range('E49').Value = '  -0.58%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.040.02'
$ws.Range('E50').Value = '  -2.85%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').Value = '2.882'
$ws.Range('E51').Value = '  -3.41%  '

$ws.Range("D2:D51").Style = "Normal"
